$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: add "Store Score" / "1. Back Bar" to the previously empty C4/D4 cells
$ws.Range("C4").Value = "Store Score"
$ws.Range("D4").Value = "1. Back Bar"

# Row 5: add "Store Score" to C5, bump the Weight (H5) from 0.25 to 0.35
$ws.Range("C5").Value = "Store Score"
$ws.Range("H5").Value = 0.35

# Row 6: add "Store Score" to C6, change Target (G6) from "N/A" text to numeric 0.25,
# and bump the Weight (H6) from 0.25 to 0.15
$ws.Range("C6").Value = "Store Score"
$ws.Range("G6").Value = 0.25
$ws.Range("H6").Value = 0.15

# Update the active selection to match the author's cursor position when saving
$ws.Range("C6").Select() | Out-Null
